$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append today's win tally as a new row under the existing header row.
$ws.Range("A2").Value = 45951
$ws.Range("A2").NumberFormat = "yyyy-mm-dd"
$ws.Range("A2").NumberFormat = "YYYY-MM-DD"

$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0
